$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in row 10 first (Sunday 29.4.18, 1200-1500, 3 hours, Insertion Sort)
$ws.Range("A10").Value = "Sunday 29.4.18"
$ws.Range("D10").Value = "Insertion Sort"
$ws.Range("B10").Value = "1200-1500"
$ws.Range("C10").Value = 3

# Fill in row 9 (Saturday 28.4.18, 2300-0100, 3 hours, Insertion Sort)
$ws.Range("A9").Value = "Saturday 28.4.18"
$ws.Range("B9").Value = "2300-0100"
$ws.Range("D9").Value = "Insertion Sort"
$ws.Range("C9").Value = 3

# Update the selected cell to A10, matching the saved view state
$ws.Range("A10").Select()
